$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder header row (row 1)
$ws.Range("A1").Value = "kitchens_1"
$ws.Range("B1").Value = "kitchens_2"
$ws.Range("C1").Value = "bedrooms_1"
$ws.Range("D1").Value = "bedrooms_2"
$ws.Range("E1").Value = "living_rooms_1"
$ws.Range("F1").Value = "living_rooms_2"

# Row 2: D2 0->1, E2 1->0
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 0

# Row 5: A5,B5,C5 from (1,0,0) to (0,0,1)
$ws.Range("A5").Value = 0
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 1

# Row 6: A6 0->1, D6 1->0
$ws.Range("A6").Value = 1
$ws.Range("D6").Value = 0

# Row 7: C7,D7,E7 from (1,0,0) to (0,0,1)
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 1
